$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F7").Value = 8
$ws.Range("F8").Value = 6
$ws.Range("F10").Value = -12
$ws.Range("F13").Value = -6
$ws.Range("F15").Value = 1
